$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for columns A and B, rows 1-32
$values = @(
    @(-0.094902791605733228, 0.094790212782953631),
    @(0.00096333064137077429, -0.0012787598773158493),
    @(0.10420937245566719, -0.1045174900675363),
    @(-0.17947492856357883, 0.17866520776028239),
    @(-0.17266520823422038, 0.17104292134561749),
    @(-0.072295378553744349, 0.072221365543691185),
    @(-0.052221366115269063, 0.052078361156768338),
    @(-0.032078361732090777, 0.032003023431348865),
    @(-0.026003023930742053, 0.025952850489696644),
    @(-0.056014786536906058, 0.055969004719926829),
    @(-0.051469005214208607, 0.051391392554332782),
    @(-0.045391393058936913, 0.045149688326154624),
    @(-0.039149688838601371, 0.039083924368398471),
    @(-0.027083924917072011, 0.027052204531104351),
    @(-0.02105220504739691, 0.021027248734728943),
    @(-0.015027249252571373, 0.015004344140423065),
    @(-0.0090043446603411681, 0.0089999994626328572),
    @(-0.085753543917171271, 0.085681237789813736),
    @(-0.076681238278197839, 0.076127275276403505),
    @(-0.067127275777965068, 0.067009963434815134),
    @(-0.058009963938929765, 0.057853776727209372),
    @(-0.093930613065877466, 0.09362321677124541),
    @(-0.08462321726287314, 0.084124538520087278),
    @(-0.04212453920678616, 0.04199999930976972),
    @(-0.092747543282722944, 0.092638406495424874),
    @(-0.086638406987884053, 0.086503053219473713),
    @(-0.080503053714116035, 0.080058270797717768),
    @(-0.074058271302134493, 0.073769466737068079),
    @(-0.063634042611809605, 0.063497609154062218),
    @(-0.0421625587146508, 0.042018910828445222),
    @(-0.027018911398076995, 0.027000782601028561),
    @(-0.00600078320543318, 0.0059999994805099988)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Update column A width (target OOXML width 15.7109375 -> nearest value the
# engine's pixel-quantized ColumnWidth can reach is 15.666666..., achieved by
# any ColumnWidth in the ~14.75-14.9167 bucket; 14.8 sits safely mid-bucket)
$ws.Columns.Item(1).ColumnWidth = 14.8
